$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert the three new "Model ROS" rows. Each Insert() pushes the existing
# row (and everything below it) down by one, which is exactly how the
# target layout was produced (a new row ends up right above each of the
# three existing blank/gap rows at 7, 12 and 17 as counted *before* any of
# these inserts happen -- i.e. row 6, then 11, then 16 after the prior
# shifts have already been applied).
# ---------------------------------------------------------------------------

$ws.Rows.Item(6).Insert()
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(16).Insert()

# ---------------------------------------------------------------------------
# Row 6 - "Model ROS" values alongside the "Wind Speed" rows
# ---------------------------------------------------------------------------
$ws.Range("A6").Value2 = "Model ROS"
$ws.Range("B6").Value2 = 6
$ws.Range("C6").Value2 = 6
$ws.Range("D6").Value2 = 2.3831376999999998
$ws.Range("E6").Value2 = 6
$ws.Range("F6").Value2 = 1.9875708000000001
$ws.Range("G6").Value2 = 5.0585839999999997
$ws.Range("I6").Value2 = 6
$ws.Range("K6").Value2 = 6.0000004999999996
$ws.Range("L6").Value2 = 6
$ws.Range("M6").Value2 = 6
$ws.Range("N6").Value2 = 6
$ws.Range("O6").Value2 = 6
$ws.Range("P6").Value2 = 6
$ws.Range("Q6").Value2 = 6
$ws.Range("R6").Value2 = 6
$ws.Range("S6").Value2 = 6.0000004999999996
$ws.Range("T6").Value2 = 6
$ws.Range("U6").Value2 = 6
$ws.Range("V6").Value2 = 6.0000004999999996

# ---------------------------------------------------------------------------
# Row 11 - "Model ROS" values alongside the "Wind Speed NE" rows
# ---------------------------------------------------------------------------
$ws.Range("A11").Value2 = "Model ROS"
$ws.Range("B11").Value2 = 0.60530806000000004
$ws.Range("C11").Value2 = 0.60672959999999998
$ws.Range("D11").Value2 = 0.40722406
$ws.Range("E11").Value2 = 0.60585829999999996
$ws.Range("F11").Value2 = 0.38986418
$ws.Range("G11").Value2 = 0.62131625000000001
$ws.Range("I11").Value2 = 0.60610549999999996
$ws.Range("K11").Value2 = 0.61175597000000004
$ws.Range("L11").Value2 = 0.52608319999999997
$ws.Range("M11").Value2 = 0.60557395000000003
$ws.Range("N11").Value2 = 0.51828969999999996
$ws.Range("O11").Value2 = 0.60281085999999995
$ws.Range("P11").Value2 = 0.60864339999999995
$ws.Range("Q11").Value2 = 0.60507619999999995
$ws.Range("R11").Value2 = 0.60420436
$ws.Range("S11").Value2 = 0.60536160000000006
$ws.Range("T11").Value2 = 0.60522880000000001
$ws.Range("U11").Value2 = 0.6053636
$ws.Range("V11").Value2 = 0.60534980000000005

# ---------------------------------------------------------------------------
# Row 16 - "Model ROS" values alongside the "Wind Speed SW" rows
# ---------------------------------------------------------------------------
$ws.Range("A16").Value2 = "Model ROS"
$ws.Range("B16").Value2 = 0.68085600000000002
$ws.Range("C16").Value2 = 0.64096624000000002
$ws.Range("E16").Value2 = 0.6803785
$ws.Range("G16").Value2 = 0.60257470000000002
$ws.Range("I16").Value2 = 0.68052029999999997
$ws.Range("K16").Value2 = 0.65760960000000002
$ws.Range("L16").Value2 = 0.62494479999999997
$ws.Range("M16").Value2 = 0.68055220000000005
$ws.Range("N16").Value2 = 0.60391843000000001
$ws.Range("O16").Value2 = 0.67636067
$ws.Range("P16").Value2 = 0.78299266000000001
$ws.Range("Q16").Value2 = 0.68069259999999998
$ws.Range("R16").Value2 = 0.75873279999999999
$ws.Range("S16").Value2 = 0.68084513999999996
$ws.Range("T16").Value2 = 0.67894129999999997
$ws.Range("U16").Value2 = 0.68056333000000002
$ws.Range("V16").Value2 = 0.67868139999999999

# ---------------------------------------------------------------------------
# Update the saved selection to match the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("S37").Select() | Out-Null
